# Daily attendance processing - 2025-11-26 22:50:09
# Reorders the comma-separated "Recorded By" values in column G so that the
# list order is reversed (e.g. "a, b" -> "b, a").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $val = $cell.Value2

    if ($null -ne $val -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Count -gt 1) {
            $reversed = @()
            for ($i = $parts.Count - 1; $i -ge 0; $i--) {
                $reversed += $parts[$i]
            }
            $newVal = $reversed -join ", "
            $cell.Value2 = $newVal
        }
    }
}
